# Experiment order generation script
# Regenerates each task-order sheet's randomized file list and reassigns
# sheet tab names so each tab's name keeps matching its (new) content.

$wb = $excel.ActiveWorkbook

# --- Sheet 1 (was GNG content, A1:B5) becomes NB content, grows to A1:B10 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "NB_TO-16515890386885605"

$ws1.Range("B2").Value = "OB-1651589036300285.csv"
$ws1.Range("B3").Value = "TB-16515890386728728.csv"
$ws1.Range("B4").Value = "ZB-match_6-1651589035779485.csv"
$ws1.Range("B5").Value = "TB-16515890370116181.csv"

# Add 5 new rows (6-10), copying formatting from the existing data rows.
$ws1.Range("A2:B5").Copy($ws1.Range("A6:B9"))
$ws1.Range("A2:B2").Copy($ws1.Range("A10:B10"))

$ws1.Range("A6").Value = 4
$ws1.Range("B6").Value = "TB-1651589036445949.csv"
$ws1.Range("A7").Value = 5
$ws1.Range("B7").Value = "ZB-match_0-16515890359616363.csv"
$ws1.Range("A8").Value = 6
$ws1.Range("B8").Value = "ZB-match_8-1651589035826366.csv"
$ws1.Range("A9").Value = 7
$ws1.Range("B9").Value = "OB-1651589036399074.csv"
$ws1.Range("A10").Value = 8
$ws1.Range("B10").Value = "OB-16515890361044407.csv"

# --- Sheet 2 (was NB content, A1:B10) becomes TOL content, shrinks to A1:B7 ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "TOL_TO-1651589038735342"

$ws2.Range("B2").Value = "MM_stims-16515890387041273.csv"
$ws2.Range("B3").Value = "ZM_stims-16515890386885605.csv"
$ws2.Range("B4").Value = "MM_stims-1651589038719717.csv"
$ws2.Range("B5").Value = "ZM_stims-16515890387041273.csv"
$ws2.Range("B6").Value = "MM_stims-1651589038735342.csv"
$ws2.Range("B7").Value = "ZM_stims-1651589038719717.csv"

$ws2.Rows("8:10").Delete()

# --- Sheet 3 (RS content, A1:B3) stays same size, values swap ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-1651589038735342"

$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

# --- Sheet 4 (was TOL content, A1:B7) becomes GNG content, shrinks to A1:B5 ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "GNG_TO-1651589038766591"

$ws4.Range("B2").Value = "go_stims-1651589038735342.csv"
$ws4.Range("B3").Value = "GNG_stims-16515890387509677.csv"
$ws4.Range("B4").Value = "go_stims-16515890387509677.csv"
$ws4.Range("B5").Value = "GNG_stims-1651589038766591.csv"

$ws4.Rows("6:7").Delete()

# --- Sheet 5 (vSAT content, A1:B5) stays same size, new values ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-16515890388290935"

$ws5.Range("B2").Value = "vSAT_stims-16515890388134685.csv"
$ws5.Range("B3").Value = "SAT_stims-1651589038766591.csv"
$ws5.Range("B4").Value = "vSAT_stims-1651589038797844.csv"
$ws5.Range("B5").Value = "SAT_stims-16515890387822154.csv"
